$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 225
$ws1.Range("F4").Value = 807
$ws1.Range("F5").Value = 237
$ws1.Range("F6").Value = 400
$ws1.Range("F7").Value = 552
$ws1.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202408/Oaqhz3fL1724668527159.jpeg"
$ws1.Range("F8").Value = 210
$ws1.Range("F11").Value = 123
$ws1.Range("F12").Value = 594
$ws1.Range("F13").Value = 77
$ws1.Range("F14").Value = 1750
$ws1.Range("F15").Value = 315
$ws1.Range("F16").Value = 2149
$ws1.Range("F17").Value = 281
$ws1.Range("F19").Value = 38
$ws1.Range("F21").Value = 129

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 221
$ws2.Range("F4").Value = 43
$ws2.Range("F7").Value = 470
$ws2.Range("F14").Value = 37

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5276
$ws3.Range("F3").Value = 304
$ws3.Range("F4").Value = 173

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5276
$ws4.Range("F4").Value = 304
$ws4.Range("F6").Value = 173
$ws4.Range("F7").Value = 225
$ws4.Range("F8").Value = 221
$ws4.Range("F9").Value = 43
$ws4.Range("F12").Value = 470
$ws4.Range("F13").Value = 807
$ws4.Range("F16").Value = 237
$ws4.Range("F17").Value = 400
$ws4.Range("F18").Value = 552
$ws4.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202408/Oaqhz3fL1724668527159.jpeg"
$ws4.Range("F19").Value = 210
$ws4.Range("F23").Value = 123
$ws4.Range("F26").Value = 594
$ws4.Range("F27").Value = 77
$ws4.Range("F29").Value = 1750
$ws4.Range("F30").Value = 315
$ws4.Range("F31").Value = 2149
$ws4.Range("F32").Value = 37
$ws4.Range("F33").Value = 281
$ws4.Range("F35").Value = 38
$ws4.Range("F38").Value = 129

$wb.Save()
